# Apply the cell-level price/volume updates for the cryptos worksheet refresh.
# All changed cells are plain text values (inline strings in the source data),
# so price cells are forced back to Text after assignment to avoid Excel
# auto-converting numeric-looking strings (e.g. "1.00" -> 1) and losing
# formatting such as trailing zeros or thousands-style separators.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.351.39"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.654.02"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.87"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.05"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.663.59"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  +9.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.335"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.121.21"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.358.90"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.04"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.640.49"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.69"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.42%  "
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.32"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.11"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0804"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.15"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  +5.12%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.74"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.51"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").Value = "  +2.22%  "
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.898"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.880"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("E40").Value = "  +3.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.58"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("E42").Value = "  +4.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "276.55"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0974"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0536"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.060.25"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("E49").Value = "  +1.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.79"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("E51").Value = "  -0.86%  "
